$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 27
$ws.Range("F5").Value = 1168
$ws.Range("F6").Value = 9152
$ws.Range("F9").Value = 7178
$ws.Range("F11").Value = 323
$ws.Range("F12").Value = 0
$ws.Range("F14").Value = 6378
$ws.Range("F15").Value = 1103
$ws.Range("F16").Value = 433
$ws.Range("F17").Value = 415
$ws.Range("F18").Value = 31
$ws.Range("F19").Value = 599
$ws.Range("F20").Value = 282
$ws.Range("F21").Value = 213
$ws.Range("F23").Value = 105
$ws.Range("F24").Value = 10293
$ws.Range("F25").Value = 1957
$ws.Range("F26").Value = 2179
$ws.Range("F28").Value = 38
$ws.Range("F29").Value = 2204
$ws.Range("F30").Value = 84
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 2130
$ws.Range("F36").Value = 314
$ws.Range("F37").Value = 1423
$ws.Range("F38").Value = 63
$ws.Range("F39").Value = 5312
$ws.Range("F40").Value = 1208
$ws.Range("G40").Value = 79
$ws.Range("F41").Value = 698
$ws.Range("F42").Value = 123
$ws.Range("F45").Value = 1079
$ws.Range("F47").Value = 1392
$ws.Range("F48").Value = 68
$ws.Range("F49").Value = 1100

# --- Sheet "演出" (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 100
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = 7

# --- Sheet "本地生活" (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 22

# --- Sheet "全部类型" (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = 27
$ws.Range("F6").Value = 17
$ws.Range("F7").Value = 1168
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 7178
$ws.Range("F16").Value = 5574
$ws.Range("F17").Value = 75
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 1103
$ws.Range("F21").Value = 433
$ws.Range("F22").Value = 415
$ws.Range("F23").Value = 599
$ws.Range("F24").Value = 282
$ws.Range("F26").Value = 0
$ws.Range("F28").Value = 10293
$ws.Range("F29").Value = 1957
$ws.Range("F30").Value = 2179
$ws.Range("F32").Value = 2204
$ws.Range("F33").Value = 0
$ws.Range("F36").Value = 74
$ws.Range("F37").Value = 2130
$ws.Range("F38").Value = 314
$ws.Range("F39").Value = 1423
$ws.Range("F40").Value = 5312
$ws.Range("G41").Value = 79
$ws.Range("F42").Value = 698
$ws.Range("F43").Value = 123
$ws.Range("F46").Value = 1079
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 68
$ws.Range("F50").Value = 1100
